{"js": "// Insert a comma after \"En esta secci\u00f3n\" so the sentence reads\n// \"En esta secci\u00f3n, se resaltar\u00e1 la diferencia...\"\nconst body = context.document.body;\n\n// Search specifically for the phrase immediately preceding the insertion\n// point (including a bit of trailing context) so we only match the\n// intended occurrence in the introductory paragraph.\nconst results = body.search(\"En esta secci\u00f3n se resaltar\u00e1\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found\");\n}\n\n// The found range spans \"En esta secci\u00f3n se resaltar\u00e1\"; insert the comma\n// right after \"En esta secci\u00f3n\" (i.e. before the 16th character) by\n// re-searching within that narrower hit for just the prefix phrase.\nconst prefixResults = results.items[0].search(\"En esta secci\u00f3n\", { matchCase: true });\nprefixResults.load(\"items\");\nawait context.sync();\n\nprefixResults.items[0].insertText(\",\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the sentence that starts the introductory paragraph. Using a Range\n# (rather than the bare Find object) means a successful Execute() collapses/\n# expands the range itself onto the matched text, so we can reposition its\n# end and insert the comma precisely after \"En esta secci\u00f3n\".\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"En esta secci\u00f3n se resaltar\u00e1\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $prefixLen = \"En esta secci\u00f3n\".Length\n    $insertPoint = $d.Range($rng.Start + $prefixLen, $rng.Start + $prefixLen)\n    $insertPoint.InsertAfter(\",\")\n}\n"}
